$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly_Update")

$ws.Range("B2").Value = -0.1
$ws.Range("C2").Value = -35.56
$ws.Range("D2").Value = 669.3133905029297
$ws.Range("B3").Value = -1.37
$ws.Range("C3").Value = 10.27
$ws.Range("D3").Value = 1227.933637446762
$ws.Range("B4").Value = 0.73
$ws.Range("C4").Value = -40.41
$ws.Range("D4").Value = 243.1582970761009
$ws.Range("B5").Value = -3.59
$ws.Range("C5").Value = -22.77
$ws.Range("D5").Value = 197.5661178723541
$ws.Range("B6").Value = 3.96
$ws.Range("C6").Value = 141.47
$ws.Range("D6").Value = 1355.89027344652
$ws.Range("B7").Value = 3.38
$ws.Range("C7").Value = -26.42
$ws.Range("D7").Value = 532.9533531413073
$ws.Range("B8").Value = 4.71
$ws.Range("C8").Value = -31.26
$ws.Range("D8").Value = 316.1528864099996
$ws.Range("B9").Value = 9.050000000000001
$ws.Range("C9").Value = -65.56999999999999
$ws.Range("D9").Value = 415.6500091552734
$ws.Range("B10").Value = 3.49
$ws.Range("C10").Value = -1.44
$ws.Range("D10").Value = 596.1
$ws.Range("B11").Value = -3.48
$ws.Range("C11").Value = 1.01
$ws.Range("D11").Value = 516.3633672706604
$ws.Range("B12").Value = -1.72
$ws.Range("C12").Value = 47.91
$ws.Range("D12").Value = 883.8513015289307
$ws.Range("B13").Value = -1.4
$ws.Range("C13").Value = -8.199999999999999
$ws.Range("D13").Value = 1362.895983432379
$ws.Range("B14").Value = 3.31
$ws.Range("C14").Value = 33.14
$ws.Range("D14").Value = 931.1642939068718
$ws.Range("B15").Value = -4.96
$ws.Range("C15").Value = -32.72
$ws.Range("D15").Value = 466.8702122518299
$ws.Range("B16").Value = 3.97
$ws.Range("C16").Value = 27.61
$ws.Range("D16").Value = 572.299296390152
$ws.Range("B17").Value = 3.05
$ws.Range("C17").Value = -71.89
$ws.Range("D17").Value = 107.1592958419346
$ws.Range("B18").Value = 0.83
$ws.Range("C18").Value = 3.95
$ws.Range("D18").Value = 10395.321715674
$ws.Range("C19").Value = 1333.465581704775
$ws.Range("D19").Value = 922.7250468660891
